$wb = $excel.ActiveWorkbook

# 1. Update status text "Ready for handoff" -> "In Translation" everywhere it appears
#    (shared string used across Overview + language sheets).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value()
        if ("Ready for handoff" -eq $v) {
            $cell.Value = "In Translation"
        }
    }
}

# 2. Narrow the "Status"-related columns that used to be sized for
#    "Ready for handoff" now that the text is shorter.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
